$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Insert 3 new product rows into the R088 list (values stored as text, same
# as every other cell in the sheet), keeping the existing thin-border style.
#
# Helper pattern per row:
#   1. Insert a blank row at the target position (existing rows shift down).
#   2. Temporarily force the row to Text format ("@") so digit-only values
#      (item codes / sequence numbers) are written as text instead of being
#      auto-coerced to numbers.
#   3. Write the six field values.
#   4. Copy the formatting (border style) from a neighboring data row and
#      paste-special just the formats back onto the new row, so the new row
#      ends up with exactly the same cell style as the rest of the table.
# ---------------------------------------------------------------------------

# New row: "G/DAY PC BB.MONSTER" -> before existing "BEAUTY TOTE BAG" row (old row 9)
$ws.Rows.Item(9).Insert()
$ws.Range("A9:F9").NumberFormat = "@"
$ws.Range("A9").Value = "20140810"
$ws.Range("B9").Value = "G/DAY PC BB.MONSTER"
$ws.Range("C9").Value = "R088"
$ws.Range("D9").Value = "1"
$ws.Range("E9").Value = "13"
$ws.Range("F9").Value = "RT"
$ws.Range("A8:F8").Copy()
$ws.Range("A9:F9").PasteSpecial(-4122)

# New row: "SUNSLK BB.MONSTER PC" -> after "BEAUTY TOTE BAG" (now row 10),
# before "KIN PHOTOCARD JKT48" (now row 11)
$ws.Rows.Item(11).Insert()
$ws.Range("A11:F11").NumberFormat = "@"
$ws.Range("A11").Value = "20140270"
$ws.Range("B11").Value = "SUNSLK BB.MONSTER PC"
$ws.Range("C11").Value = "R088"
$ws.Range("D11").Value = "2"
$ws.Range("E11").Value = "2"
$ws.Range("F11").Value = "TG"
$ws.Range("A10:F10").Copy()
$ws.Range("A11:F11").PasteSpecial(-4122)

# New row: "FIESTA STANDEE TREAS" -> appended after "KIN PTCARD JKT48 SRS"
# (now row 14), becoming the new last row 15
$ws.Range("A15:F15").NumberFormat = "@"
$ws.Range("A15").Value = "20141232"
$ws.Range("B15").Value = "FIESTA STANDEE TREAS"
$ws.Range("C15").Value = "R088"
$ws.Range("D15").Value = "4"
$ws.Range("E15").Value = "4"
$ws.Range("F15").Value = "RT"
$ws.Range("A14:F14").Copy()
$ws.Range("A15:F15").PasteSpecial(-4122)

$excel.CutCopyMode = 0
